# The workbook has two result tables on "Лист1" (sheet 1):
#   - one headed "Неявная схема"              (merged A1:C1)
#   - one headed "Схема Кранка— Никольсона"   (merged A7:C7)
# This edit swaps the two headings (so the Crank-Nicolson block now comes
# first, at A1, and the implicit-scheme block is second, at A7), and widens
# the second (now "Неявная схема") header band from 3 columns to 6 columns
# (A7:C7 plus a new D7:F7 merge) with the same centered formatting.
# The active selection is also left on D12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the two section headings -----------------------------------
$topHeading = $ws.Range("A1").Value()
$bottomHeading = $ws.Range("A7").Value()

$ws.Range("A1").Value = $bottomHeading
$ws.Range("A7").Value = $topHeading

# --- Extend the (now) "Неявная схема" header band to D7:F7 -----------
# Match the centered formatting already used on A7:C7, then merge.
$ws.Range("D7:F7").HorizontalAlignment = $ws.Range("A7").HorizontalAlignment
$null = $ws.Range("D7:F7").Merge()

# --- Update the selected cell -----------------------------------------
$null = $ws.Range("D12").Select()
